# The data table on Sheet1 (Hortaliza, Feria Lagunitas de Puerto Montt - Coliflor)
# gets one new weekly record inserted at row 269. Inserting a row there shifts
# every following record down by one (row 270 becomes the old row 269, etc.),
# which also naturally produces the new trailing row 383 (a duplicate of the
# old last row, 382) and bumps the sheet dimension from A1:R382 to A1:R383.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before row 269; Excel shifts rows 269:382 down to 270:383.
$ws.Rows("269:269").Insert()

# Populate the newly inserted row 269 with the new weekly record.
$ws.Cells.Item(269, 1).Value = 4
$ws.Cells.Item(269, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(269, 3).Value = 'Los Lagos'
$ws.Cells.Item(269, 4).Value = 44784
$ws.Cells.Item(269, 5).Value = 10
$ws.Cells.Item(269, 6).Value = 100112008
$ws.Cells.Item(269, 7).Value = 'Coliflor'
$ws.Cells.Item(269, 8).Value = 'Sin especificar'
$ws.Cells.Item(269, 9).Value = 'Primera'
$ws.Cells.Item(269, 10).Value = 500
$ws.Cells.Item(269, 11).Value = 1700
$ws.Cells.Item(269, 12).Value = 1700
$ws.Cells.Item(269, 13).Value = 1700
$ws.Cells.Item(269, 14).Value = '$/unidad'
$ws.Cells.Item(269, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(269, 16).Value = 1700
$ws.Cells.Item(269, 17).Value = 1
$ws.Cells.Item(269, 18).Value = 'Hortaliza'
